# Update "想去人数" (F) and "最低票价" (G) figures for a handful of
# conventions on both the "展览" and "全部类型" worksheets, matching the
# freshly generated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# Row -> column -> new value, applied identically to both sheets.
$updates = @(
    @{ Row = 3;  Col = "G"; Value = 50 }
    @{ Row = 4;  Col = "G"; Value = 65 }
    @{ Row = 8;  Col = "F"; Value = 33 }
    @{ Row = 12; Col = "F"; Value = 1138 }
    @{ Row = 18; Col = "F"; Value = 160 }
    @{ Row = 22; Col = "F"; Value = 280 }
    @{ Row = 28; Col = "F"; Value = 633 }
    @{ Row = 30; Col = "F"; Value = 92 }
    @{ Row = 31; Col = "F"; Value = 3973 }
    @{ Row = 33; Col = "F"; Value = 465 }
    @{ Row = 35; Col = "F"; Value = 1015 }
    @{ Row = 36; Col = "F"; Value = 106 }
    @{ Row = 39; Col = "F"; Value = 111 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $cellAddr = "$($u.Col)$($u.Row)"
        $ws.Range($cellAddr).Value = $u.Value
    }
}
